# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The worker table (rows 16-19) is reordered: the previously-first worker
# ("ZOBEIDA CASTRO POLO") moves to the bottom of the block, and the three
# other workers shift up to take the top three rows. Each worker's own
# Periodo/Valor Mora/Salario Basico values travel with them (no values are
# recombined), so this is a pure row reorder of the table body.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 16: YANILETH ESTHER MEDINA BETANCOURT (was row 17)
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1047434072"
$ws.Range("D16").Value = "YANILETH ESTHER MEDINA BETANCOURT"
$ws.Range("E16").Value = "2202"
$ws.Range("F16").Value = 20000
$ws.Range("G16").Value = 1300000

# New row 17: CARLOS ENRIQUE QUINTANA SLAGADO (was row 18)
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1082241049"
$ws.Range("D17").Value = "CARLOS ENRIQUE QUINTANA SLAGADO"
$ws.Range("E17").Value = "2202"
$ws.Range("F17").Value = 18666
$ws.Range("G17").Value = 2000000

# New row 18: JOSE DANIEL VARGAS PAJARO (was row 19)
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1007739316"
$ws.Range("D18").Value = "JOSE DANIEL VARGAS PAJARO"
$ws.Range("E18").Value = "2202"
$ws.Range("F18").Value = 20000
$ws.Range("G18").Value = 877803

# New row 19: ZOBEIDA CASTRO POLO (was row 16) -- moved to the bottom
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "52406183"
$ws.Range("D19").Value = "ZOBEIDA CASTRO POLO"
$ws.Range("E19").Value = "2108"
$ws.Range("F19").Value = 36341
$ws.Range("G19").Value = 908528

# Row 20 (DANNES ELENA CASTILLO MORENO) is unchanged.
